$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A34").Value = "팅팅 씨, 어디에 가요?"
$ws.Range("C34").Value = "팅팅씨어디에가요.mp3"

$ws.Range("A35").Value = "식당에 가요. 식당에서 친구하고 밥을 먹어요."
$ws.Range("C35").Value = "식당에가요식당에서친구하고밥을먹어요.mp3"

$ws.Range("A36").Value = "아, 그래요?"
$ws.Range("C36").Value = "아그래요.mp3"

$ws.Range("A37").Value = "케빈 씨는 오늘 뭐 해요?"
$ws.Range("C37").Value = "케빈씨는오늘뭐해요.mp3"

$ws.Range("A38").Value = "저는 집에서 쉬어요."
$ws.Range("C38").Value = "저는집에서쉬어요.mp3"
